$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.233.29"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3
$ws.Range("D3").Value = "2.360.50"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.36%  "

# Row 5
$ws.Range("D5").Value = "'550.64"
$ws.Range("E5").Value = "  +0.90%  "

# Row 6
$ws.Range("D6").Value = "'133.18"
$ws.Range("E6").Value = "  -1.90%  "

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.107"
$ws.Range("E9").Value = "  +4.58%  "

# Row 10
$ws.Range("D10").Value = "'5.68"
$ws.Range("E10").Value = "  +4.95%  "

# Row 11
$ws.Range("E11").Value = "  -1.33%  "

# Row 12
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  -0.94%  "

# Row 13
$ws.Range("D13").Value = "'24.20"
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("D14").Value = "2.772.83"
$ws.Range("E14").Value = "  -0.04%  "

# Row 15
$ws.Range("D15").Value = "57.962.13"
$ws.Range("E15").Value = "  -0.29%  "

# Row 16
$ws.Range("E16").Value = "  +2.09%  "

# Row 17
$ws.Range("D17").Value = "2.381.33"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18
$ws.Range("D18").Value = "'11.04"
$ws.Range("E18").Value = "  +3.60%  "

# Row 19
$ws.Range("D19").Value = "'4.33"
$ws.Range("E19").Value = "  +1.86%  "

# Row 20
$ws.Range("D20").Value = "'331.69"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").Value = "'6.92"
$ws.Range("E21").Value = "  +3.18%  "

# Row 22
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "'63.74"
$ws.Range("E23").Value = "  +2.31%  "

# Row 24
$ws.Range("D24").Value = "'0.170"
$ws.Range("E24").Value = "  +1.21%  "

# Row 25
$ws.Range("E25").Value = "  +0.40%  "

# Row 26
$ws.Range("D26").Value = "'8.30"
$ws.Range("E26").Value = "  -2.90%  "

# Row 27
$ws.Range("D27").Value = "'1.33"
$ws.Range("E27").Value = "  -7.21%  "

# Row 28
$ws.Range("D28").Value = "'1.77"
$ws.Range("E28").Value = "  -0.29%  "

# Row 29
$ws.Range("D29").Value = "'170.49"
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0742"
$ws.Range("E30").Value = "  +0.54%  "

# Row 31
$ws.Range("D31").Value = "'6.18"
$ws.Range("E31").Value = "  -0.08%  "

# Row 32
$ws.Range("D32").Value = "'18.43"
$ws.Range("E32").Value = "  -0.68%  "

# Row 33
$ws.Range("E33").Value = "  -3.57%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  -0.42%  "

# Row 36
$ws.Range("D36").Value = "'4.17"
$ws.Range("E36").Value = "  -0.62%  "

# Row 37
$ws.Range("E37").Value = "  -1.48%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'40.31"
$ws.Range("E38").Value = "  +2.88%  "

# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.427"
$ws.Range("E39").Value = "  +13.31%  "

# Row 40
$ws.Range("D40").Value = "'1.60"
$ws.Range("E40").Value = "  -1.80%  "

# Row 41
$ws.Range("D41").Value = "'3.68"
$ws.Range("E41").Value = "  +1.23%  "

# Row 42
$ws.Range("D42").Value = "'140.90"
$ws.Range("E42").Value = "  -4.61%  "

# Row 43
$ws.Range("D43").Value = "'288.31"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").Value = "'0.0959"
$ws.Range("E44").Value = "  +2.01%  "

# Row 45
$ws.Range("D45").Value = "'0.0516"
$ws.Range("E45").Value = "  +1.68%  "

# Row 46
$ws.Range("E46").Value = "  +0.34%  "

# Row 47
$ws.Range("D47").Value = "'0.399"
$ws.Range("E47").Value = "  +3.53%  "

# Row 48
$ws.Range("D48").Value = "'18.63"
$ws.Range("E48").Value = "  -2.10%  "

# Row 49
$ws.Range("D49").Value = "'0.0223"
$ws.Range("E49").Value = "  +2.43%  "

# Row 50
$ws.Range("E50").Value = "  -0.52%  "

# Row 51
$ws.Range("E51").Value = "  +0.10%  "
